# Adds a new "CLASSROOMS" worksheet at the end of the workbook, listing
# classroom numbers and their normal seating capacity.
# (commit: "added doubly linked list/node and updated get classroom function")

$wb = $excel.ActiveWorkbook

# Append a brand-new worksheet after the last existing sheet so it becomes
# sheet #9 / the last tab, and make it the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CLASSROOMS"

# Header row
$ws.Cells.Item(1, 1).Value = "Classroom #"
$ws.Cells.Item(1, 2).Value = "Normal Capacity"

# Classroom # / Normal Capacity data
$rooms = @(
    @("11-533", 36),
    @("11-534", 36),
    @("11-560", 24),
    @("11-562", 24),
    @("11-564", 24),
    @("11-458", 40),
    @("11-430 ", 30),
    @("11-320", 30),
    @("11-532 Computer Lab", 30)
)

$r = 2
foreach ($room in $rooms) {
    $ws.Cells.Item($r, 1).Value = $room[0]
    $ws.Cells.Item($r, 2).Value = $room[1]
    $r++
}

# Column widths to comfortably fit the room names / capacity labels.
$ws.Columns.Item(1).ColumnWidth = 33
$ws.Columns.Item(2).ColumnWidth = 39.17

# Leave the selection on the last entered cell, like the source workbook.
$ws.Range("A10").Select() | Out-Null
